$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 204; this shifts the existing row 204
# (and everything below it, through the former row 302) down by one,
# to rows 205-303, growing the used range from A1:R302 to A1:R303.
$ws.Rows.Item(204).Insert()

# Populate the newly inserted row 204 with the new weekly price record.
$ws.Range("A204").Value = 3
$ws.Range("B204").Value = "Femacal de La Calera"
$ws.Range("C204").Value = "Coquimbo"
$ws.Range("D204").Value = 44609
$ws.Range("E204").Value = 5
$ws.Range("F204").Value = 100112040
$ws.Range("G204").Value = "Cilantro"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 155
$ws.Range("K204").Value = 5000
$ws.Range("L204").Value = 5500
$ws.Range("M204").Value = 5242
$ws.Range("N204").Value = "$/docena de atados (3 kilos)"
$ws.Range("O204").Value = "Provincia de Quillota"
$ws.Range("P204").Value = 1747
$ws.Range("Q204").Value = 3
$ws.Range("R204").Value = "Hortaliza"
